$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.292.57"
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = "'3.420.92"
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'577.14"
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = "'148.34"
$ws.Range("E6").Value = '  +1.28%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'0.485"
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("D9").Value = "'8.03"
$ws.Range("E9").Value = '  +5.17%  '
$ws.Range("D10").Value = "'0.123"
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("D11").Value = "'0.414"
$ws.Range("E11").Value = '  +3.71%  '
$ws.Range("D12").Value = "'4.009.75"
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").Value = "'28.28"
$ws.Range("E14").Value = '  -4.10%  '
$ws.Range("D15").Value = "'3.485.82"
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").Value = "'62.388.42"
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = "'6.37"
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("D19").Value = "'14.42"
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("D20").Value = "'8.92"
$ws.Range("E20").Value = '  -2.56%  '
$ws.Range("D21").Value = "'383.60"
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").Value = "'0.566"
$ws.Range("E22").Value = '  +2.00%  '
$ws.Range("D23").Value = "'74.73"
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = "'3.589.62"
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").Value = "'0.0000112"
$ws.Range("E26").Value = '  -1.58%  '
$ws.Range("D27").Value = "'0.183"
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("D28").Value = "'7.62"
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").Value = "'7.93"
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").Value = "'2.12"
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = "'1.33"
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("D34").Value = "'23.13"
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").Value = "'5.43"
$ws.Range("E35").Value = '  +3.77%  '
$ws.Range("D36").Value = "'1.63"
$ws.Range("E36").Value = '  +4.48%  '
$ws.Range("D37").Value = "'31.42"
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = "'6.89"
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("D39").Value = "'169.01"
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("D40").Value = "'3.456.52"
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").Value = "'0.0782"
$ws.Range("E41").Value = '  +4.33%  '
$ws.Range("D42").Value = "'0.782"
$ws.Range("E42").Value = '  -1.57%  '
$ws.Range("D43").Value = "'42.43"
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = "'4.37"
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").Value = "'1.68"
$ws.Range("E45").Value = '  -1.37%  '
$ws.Range("D46").Value = "'1.16"
$ws.Range("E46").Value = '  -2.43%  '
$ws.Range("D47").Value = "'2.540.94"
$ws.Range("E47").Value = '  -2.06%  '
$ws.Range("D48").Value = "'6.90"
$ws.Range("E48").Value = '  +2.94%  '
$ws.Range("D49").Value = "'2.21"
$ws.Range("D50").Value = "'22.50"
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("E51").Value = '  +0.20%  '
